# Pio's ERA operativo Abril-Diciembre 2025
# Updates the computed workload-distribution figures on Sheet1 following
# the refreshed vacation/assignment data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Fernandez
$ws.Cells.Item(2, 7).Value = 0.8833922261484137
$ws.Cells.Item(2, 12).Value = 3.333333333333333
$ws.Cells.Item(2, 13).Value = 2
$ws.Cells.Item(2, 19).Value = -0.003684312584836145

# Row 3 - Gomez
$ws.Cells.Item(3, 7).Value = 2.347299343765784
$ws.Cells.Item(3, 8).Value = 39
$ws.Cells.Item(3, 15).Value = 4.879012345679012
$ws.Cells.Item(3, 19).Value = 0.3336076817558298

# Row 4 - Bravo
$ws.Cells.Item(4, 7).Value = 2.179034157832754
$ws.Cells.Item(4, 12).Value = 2.666666666666667
$ws.Cells.Item(4, 13).Value = 2
$ws.Cells.Item(4, 14).Value = 2
$ws.Cells.Item(4, 15).Value = 4.360655737704918
$ws.Cells.Item(4, 19).Value = 0.3322404371584699

# Row 5 - Iñiguez
$ws.Cells.Item(5, 7).Value = 2.305233047282527
$ws.Cells.Item(5, 8).Value = 38
$ws.Cells.Item(5, 15).Value = 4.639357429718875
$ws.Cells.Item(5, 19).Value = 0.3255689424364123

# Row 6 - Breinbauer
$ws.Cells.Item(6, 7).Value = 2.305233047282527
$ws.Cells.Item(6, 8).Value = 38
$ws.Cells.Item(6, 15).Value = 4.676923076923077
$ws.Cells.Item(6, 19).Value = 0.3282051282051282

# Row 7 - Arredondo
$ws.Cells.Item(7, 7).Value = 7.580346626283051
$ws.Cells.Item(7, 8).Value = 35
$ws.Cells.Item(7, 11).Value = 15
$ws.Cells.Item(7, 15).Value = 4.32520325203252
$ws.Cells.Item(7, 18).Value = 1.853658536585366
$ws.Cells.Item(7, 19).Value = 0.4943089430894309

# Row 8 - Carrasco
$ws.Cells.Item(8, 7).Value = 7.874810701665854
$ws.Cells.Item(8, 9).Value = 35
$ws.Cells.Item(8, 11).Value = 15
$ws.Cells.Item(8, 16).Value = 4.273092369477911
$ws.Cells.Item(8, 18).Value = 1.831325301204819
$ws.Cells.Item(8, 19).Value = 0.4883534136546185

# Row 9 - Culaciati
$ws.Cells.Item(9, 7).Value = 8.211341073531914
$ws.Cells.Item(9, 9).Value = 35
$ws.Cells.Item(9, 11).Value = 16
$ws.Cells.Item(9, 16).Value = 4.32520325203252
$ws.Cells.Item(9, 18).Value = 1.977235772357723
$ws.Cells.Item(9, 19).Value = 0.4943089430894309

# Row 10 - Contreras
$ws.Cells.Item(10, 7).Value = 8.699310112737702
$ws.Cells.Item(10, 9).Value = 38
$ws.Cells.Item(10, 11).Value = 17
$ws.Cells.Item(10, 16).Value = 4.695934959349593
$ws.Cells.Item(10, 18).Value = 2.100813008130081
$ws.Cells.Item(10, 19).Value = 0.4943089430894309

# Row 11 - Cisternas
$ws.Cells.Item(11, 7).Value = 8.749789668517613
$ws.Cells.Item(11, 9).Value = 39
$ws.Cells.Item(11, 11).Value = 17
$ws.Cells.Item(11, 16).Value = 4.542528735632184
$ws.Cells.Item(11, 18).Value = 1.98007662835249
$ws.Cells.Item(11, 19).Value = 0.4659003831417624

# Row 12 - Pio
$ws.Cells.Item(12, 7).Value = 7.954736664984045
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 6
$ws.Cells.Item(12, 11).Value = 17
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 17).Value = 0.7475409836065573
$ws.Cells.Item(12, 18).Value = 2.118032786885246
$ws.Cells.Item(12, 19).Value = 0.8721311475409832

# Row 13 - Alvo
$ws.Cells.Item(13, 7).Value = 7.954736664984045
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 10).Value = 6
$ws.Cells.Item(13, 11).Value = 17
$ws.Cells.Item(13, 12).Value = 6.999999999999997
$ws.Cells.Item(13, 13).Value = 2
$ws.Cells.Item(13, 14).Value = 4
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 17).Value = 0.7384615384615385
$ws.Cells.Item(13, 18).Value = 2.092307692307692
$ws.Cells.Item(13, 19).Value = 0.8615384615384613

# Row 14 - Boettiger
$ws.Cells.Item(14, 7).Value = 7.954736664984045
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 10).Value = 6
$ws.Cells.Item(14, 11).Value = 17
$ws.Cells.Item(14, 12).Value = 6.999999999999997
$ws.Cells.Item(14, 14).Value = 2
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 17).Value = 0.7414634146341462
$ws.Cells.Item(14, 18).Value = 2.100813008130081
$ws.Cells.Item(14, 19).Value = 0.8650406504065037

# Row 15 - Loch
$ws.Cells.Item(15, 7).Value = 8.379606259464948
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 6
$ws.Cells.Item(15, 11).Value = 18
$ws.Cells.Item(15, 12).Value = 7.33333333333333
$ws.Cells.Item(15, 13).Value = 2
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 17).Value = 0.7384615384615385
$ws.Cells.Item(15, 18).Value = 2.215384615384616
$ws.Cells.Item(15, 19).Value = 0.9025641025641022

# Row 16 - Rubio
$ws.Cells.Item(16, 5).Value = 14
$ws.Cells.Item(16, 7).Value = 8.442705704189834
$ws.Cells.Item(16, 11).Value = 18
$ws.Cells.Item(16, 12).Value = 7.33333333333333
$ws.Cells.Item(16, 14).Value = 2
$ws.Cells.Item(16, 17).Value = 0.8153256704980842
$ws.Cells.Item(16, 18).Value = 2.096551724137931
$ws.Cells.Item(16, 19).Value = 0.8541507024265641

# Row 17 - Recluta1
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 7).Value = 8.177688036345311
$ws.Cells.Item(17, 11).Value = 18
$ws.Cells.Item(17, 12).Value = 6.333333333333331
$ws.Cells.Item(17, 13).Value = 2
$ws.Cells.Item(17, 14).Value = 3
$ws.Cells.Item(17, 17).Value = 0.8184615384615385
$ws.Cells.Item(17, 18).Value = 2.104615384615385
$ws.Cells.Item(17, 19).Value = 0.7405128205128203
